$wb = $excel.ActiveWorkbook

# --- Rename "Requested quantity" headers to their series-specific names ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after "Monthly Trend" ---
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
  @(45025.99999999999, 6, -1.649507688058842, 13.99965051313022),
  @(45032.99999999999, 6, -1.780714193395614, 14.60923390324007),
  @(45039.99999999999, 6, -1.995813850873737, 13.97745743013591),
  @(45046.99999999999, 7, -1.551827469590339, 14.89528231258715),
  @(45053.99999999999, 7, -1.692894882863353, 15.0525752968645),
  @(45060.99999999999, 7, -1.261046901424295, 15.04968732908005),
  @(45067.99999999999, 7, -1.411184639654876, 14.44729181103635),
  @(45074.99999999999, 7, -0.7989722609515357, 14.77591894423036),
  @(45081.99999999999, 7, -0.698805193601624, 15.31504238057353),
  @(45088.99999999999, 7, -0.6877502982667604, 16.24190265202412),
  @(45095.99999999999, 7, -0.4737069500067956, 15.37337429706278),
  @(45102.99999999999, 8, -0.5199865175914298, 16.03781243281203),
  @(45116.99999999999, 8, -0.3141843027113468, 15.96478720092163),
  @(45123.99999999999, 8, 0.4459873659112629, 15.95439128290744),
  @(45151.99999999999, 9, 0.177371475552188, 16.70851377048415),
  @(45200.99999999999, 9, 1.288061133448809, 17.41859415171614),
  @(45214.99999999999, 10, 1.451536065663199, 17.98823542975003),
  @(45221.99999999999, 10, 2.133049588171711, 18.10456024100484),
  @(45235.99999999999, 10, 1.487264407788123, 18.44721833641411),
  @(45242.99999999999, 10, 1.634838125710373, 18.39442604702621),
  @(45249.99999999999, 10, 2.161834844047659, 17.70450908595315),
  @(45340.99999999999, 12, 3.42481549216728, 19.9819412363247),
  @(45354.99999999999, 12, 4.804002364803916, 19.95650052529909),
  @(45361.99999999999, 12, 4.403747675589249, 20.53123874242164),
  @(45368.99999999999, 13, 4.445774325996467, 20.94992268330127),
  @(45375.99999999999, 13, 4.154318263616167, 20.96060593778549),
  @(45382.99999999999, 13, 4.798887168987305, 21.01249253383548),
  @(45389.99999999999, 13, 4.791032322384114, 20.75382066970968),
  @(45396.99999999999, 13, 4.837872751447502, 20.87880341099856),
  @(45403.99999999999, 13, 4.869497089119434, 21.3909615583924),
  @(45410.99999999999, 13, 4.501792000883549, 21.91896171325545),
  @(45417.99999999999, 14, 4.970226807775261, 21.87481191893497)
)

$r = 2
foreach ($row in $data) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

$lastRow = $r - 1

# --- Formatting: mirror the header/date styling used on the other sheets ---
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A" + $lastRow).PasteSpecial(-4122)
